# Prepare for experiment 2025-10-14/0000
#
# Log the new "models" config run that was produced for the
# 2025-10-14/b/0000 experiment: a GRU whose rnn_input_size mirrors the
# embedding dim (like the other post-2025-09-09 GRU/RNN rows), paired
# with a smaller single-layer Identity readout network.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("models")

$newRow = 8

$ws.Cells.Item($newRow, 1).Value  = "2025-10-14 18:07:57"   # timestamp
$ws.Cells.Item($newRow, 2).Value  = "models"                # config_kind
$ws.Cells.Item($newRow, 3).Value  = "2025-10-14/b/0000"      # config_id
# note: kept as an explicit empty string, like every other logged row.
$ws.Cells.Item($newRow, 4).Value  = "'"
$ws.Cells.Item($newRow, 5).Value  = "models.networks.FCN"    # input_network
# input_network_layer_sizes: left blank for this run.
$ws.Cells.Item($newRow, 6).Value  = "'"
$ws.Cells.Item($newRow, 7).Value  = "[CallableConfig(path='torch.nn.modules.activation.ReLU', args_cfg=ReLUConfig(inplace=False), kind='class', recovery_mode='call', locked=False, if_recover_while_locked='print')]"
$ws.Cells.Item($newRow, 8).Value  = "[None]"                 # input_network_dropouts
$ws.Cells.Item($newRow, 9).Value  = "torch.nn.modules.rnn.GRU"  # rnn_type
$ws.Cells.Item($newRow, 10).Value = "embedding_dim___"       # rnn_input_size
$ws.Cells.Item($newRow, 11).Value = 20                       # rnn_hidden_size
# rnn_nonlinearity: GRU has none, like the other GRU rows.
$ws.Cells.Item($newRow, 12).Value = "'"
$ws.Cells.Item($newRow, 13).Value = "models.networks.FCN"    # readout_network
$ws.Cells.Item($newRow, 14).Value = "[20, 2]"                # readout_network_layer_sizes
$ws.Cells.Item($newRow, 15).Value = "[CallableConfig(path='torch.nn.modules.linear.Identity', args_cfg=IdentityConfig(), kind='class', recovery_mode='call', locked=False, if_recover_while_locked='print')]"
$ws.Cells.Item($newRow, 16).Value = "[None]"                 # readout_network_dropouts

# The "training" sheet's logged rows were also re-stamped with an
# explicit (visually-default) cell style when the workbook was prepared
# for this experiment.
$wsTraining = $wb.Worksheets.Item("training")
$wsTraining.Range("A2:N26").Style = "Normal"
